$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet
$ws.Name = "Sheet1"

# Row 2 (Arsenal) updates
$ws.Range("E2").Value = 27
$ws.Range("F2").Value = 297
$ws.Range("G2").Value = 2430
$ws.Range("H2").Value = 27
$ws.Range("I2").Value = 48
$ws.Range("J2").Value = 36
$ws.Range("K2").Value = 84
$ws.Range("L2").Value = 45
$ws.Range("O2").Value = 37
$ws.Range("Q2").Value = 1.78
$ws.Range("R2").Value = 1.33
$ws.Range("S2").Value = 3.11
$ws.Range("T2").Value = 1.67
$ws.Range("U2").Value = 3

# Row 21 (Wolves) updates
$ws.Range("B21").Value = 29
$ws.Range("E21").Value = 27
$ws.Range("F21").Value = 297
$ws.Range("G21").Value = 2430
$ws.Range("H21").Value = 27
$ws.Range("I21").Value = 17
$ws.Range("J21").Value = 11
$ws.Range("K21").Value = 28
$ws.Range("L21").Value = 15
$ws.Range("O21").Value = 57
$ws.Range("Q21").Value = 0.63
$ws.Range("R21").Value = 0.41
$ws.Range("S21").Value = 1.04
$ws.Range("T21").Value = 0.5600000000000001
$ws.Range("U21").Value = 0.96
